$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-09-06 21:06:35"

$wsZhCn.Range("H3").Value = "2016-09-06 21:06:30"
$wsZhCn.Range("K3").Value = "2016-09-06 21:06:47"

$wsDeDe.Range("H3").Value = "2016-09-06 21:06:35"
$wsDeDe.Range("K3").Value = "2016-09-06 21:06:55"
